# Outstandings.xlsx edit script
# Applies:
#  - Sheet1 (Purchase 22-23): add a new invoice line for "Namrata Rubber Product Pvt Ltd"
#    (group Sr.No 1), remove the stray "1693" line in the Aquachemitech group (merging the
#    vendor note onto the remaining line), and remove the "Vissu Virgincar & Sons" group (Sr.No 6).
#  - Sheet2 (Sale 22-23): adjust two outstanding-balance formulas.
#  - Selection / active tab bookkeeping to match what the workbook looked like after the edit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Purchase 22-23")
$ws2 = $wb.Worksheets.Item("Sale 22-23")

# ---------------------------------------------------------------------------
# Sheet1: insert a new row right after row 3 for a second "Namrata Rubber
# Product Pvt Ltd" invoice, copying row 3's look & feel.
# ---------------------------------------------------------------------------
$ws1.Rows("4:4").Insert()
$ws1.Range("A3:F3").Copy()
$ws1.Range("A4:F4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws1.Range("B4").Value = 45355
$ws1.Range("C4").Value = "114/23-24"
$ws1.Range("D4").Value = "Namrata Rubber Product Pvt Ltd"
$ws1.Range("E4").Value = 47466

# Row 3 no longer carries its own subtotal; the subtotal now lives on row 4,
# covering both invoices of the group.
$ws1.Range("F3").ClearContents()
$ws1.Range("F4").Formula = "=E3+E4"

# ---------------------------------------------------------------------------
# Sheet1: remove the extra "1693" Aquachemitech line (old row 16, now row 17
# after the insert above). The surviving "1530" line (row 16) becomes a
# standalone entry: its amount no longer needs the red highlight, and its
# subtotal formula now just refers to itself.
# ---------------------------------------------------------------------------
$ws1.Range("E14").Copy()
$ws1.Range("E16").PasteSpecial(-4122)     # xlPasteFormats (drop the red highlight)
$excel.CutCopyMode = 0

$ws1.Range("F16").Formula = "=E16"
$ws1.Rows("17:17").Delete()

# ---------------------------------------------------------------------------
# Sheet1: remove the whole "Vissu Virgincar & Sons" group (Sr.No 6, old row
# 18, now row 18 still since it sat below the deletion above).
# ---------------------------------------------------------------------------
$ws1.Rows("18:18").Delete()

# ---------------------------------------------------------------------------
# Sheet2: tighten two outstanding-balance formulas.
# ---------------------------------------------------------------------------
$ws2.Range("F13").Formula = "=E13-175496-500000-800000-200000"
$ws2.Range("F20").Formula = "=E20-20000-17000"

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping: Purchase 22-23 becomes the visible,
# active tab (with F3 selected); Sale 22-23 keeps F14 selected but is no
# longer the active tab.
# ---------------------------------------------------------------------------
$ws2.Range("F14").Select()
$ws1.Select()
$ws1.Range("F3").Select()
